$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "2020-01-01,2020-12-31"
$ws.Range("D8").Value = "2020-01-01,2020-12-31"
$ws.Range("B9").Value = "Pick a time (24 hrs)"

$ws.Range("B10").Select()
